# Add 2022-Q3 data:
#  - duplicate the "2022-Q1" sheet (most recent at the time) to create the new
#    "2022-Q3" sheet, positioned right after "总计" / right before "2022-Q1",
#    then overwrite its financial figures with the new quarter's numbers.
#  - insert a new summary row in "总计" for 2022-Q3, pushing the existing
#    quarters down by one row.

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $addr, $val) {
    # Force the value to be written back as text (matching the existing
    # inline/shared-string cells in this column) instead of letting Excel
    # auto-coerce a numeric-looking string into a number.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Create the "2022-Q3" worksheet from a copy of "2022-Q1"
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$q1Sheet.Copy($q1Sheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

Set-TextValue $q3Sheet "D2" "47.57"
Set-TextValue $q3Sheet "E2" "73.75"
Set-TextValue $q3Sheet "F2" "2.06"
Set-TextValue $q3Sheet "G2" "0.9799"

Set-TextValue $q3Sheet "D3" "0.27"
Set-TextValue $q3Sheet "E3" "73.75"
Set-TextValue $q3Sheet "F3" "2.06"
Set-TextValue $q3Sheet "G3" "0.0056"

# ---------------------------------------------------------------------
# 2. Insert the matching summary row into "总计"
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Copy the style of the (now-shifted) neighbouring row's index cell so the
# new A2 keeps the same bold/centered/bordered look as the rest of column A.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").Style = "Normal"

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.99

# Renumber the (0-based) index column for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
